# Updated symbol list on Thu Feb 16 14:57:57 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns of the
# crypto ticker sheet with newly scraped values. The source sheet stores
# these figures as literal text (e.g. "317.50", "4.38%") rather than native
# Excel numbers, so each value is written with a leading apostrophe to force
# Excel to keep it as text instead of auto-converting it to a Number /
# Percentage cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $ws.Range($address).Value = "'" + $value
}

Set-TextValue "D2"  "318.43"
Set-TextValue "E2"  "5.02%"

Set-TextValue "D3"  "48.15"
Set-TextValue "E3"  "11.69%"

Set-TextValue "D4"  "5.288"
Set-TextValue "E4"  "3.92%"

Set-TextValue "D5"  "0.07940"
Set-TextValue "E5"  "3.46%"

Set-TextValue "D6"  "4.591"
Set-TextValue "E6"  "4.10%"

Set-TextValue "D7"  "1.349"
Set-TextValue "E7"  "32.92%"

Set-TextValue "D8"  "1.633"
Set-TextValue "E8"  "1.20%"

Set-TextValue "D9"  "0.1293"
Set-TextValue "E9"  "3.85%"

Set-TextValue "D10" "0.1952"
Set-TextValue "E10" "5.39%"

Set-TextValue "D11" "0.09430"
Set-TextValue "E11" "2.82%"

Set-TextValue "D12" "0.04609"
Set-TextValue "E12" "10.56%"

Set-TextValue "E13" "0.06%"

Set-TextValue "D14" "0.001319"
Set-TextValue "E14" "1.84%"

Set-TextValue "D15" "0.04161"
Set-TextValue "E15" "-0.32%"

Set-TextValue "E16" "2.24%"

Set-TextValue "E17" "0.15%"

Set-TextValue "D18" "2.423"
Set-TextValue "E18" "2.88%"

Set-TextValue "D19" "0.3461"
Set-TextValue "E19" "3.22%"

Set-TextValue "D20" "8.117"
Set-TextValue "E20" "-3.58%"

Set-TextValue "D21" "0.1383"
Set-TextValue "E21" "-0.98%"

Set-TextValue "D22" "0.3099"
Set-TextValue "E22" "-2.93%"

Set-TextValue "D23" "0.001319"
Set-TextValue "E23" "2.52%"

Set-TextValue "D24" "0.004251"
Set-TextValue "E24" "-5.48%"

Set-TextValue "D25" "0.0001351"
Set-TextValue "E25" "-0.10%"

Set-TextValue "D26" "0.0003544"

Set-TextValue "D38" "0.02664"
Set-TextValue "E38" "8.32%"

Set-TextValue "D39" "0.05701"
Set-TextValue "E39" "8.06%"

Set-TextValue "D40" "0.01076"
Set-TextValue "E40" "79.93%"

Set-TextValue "D41" "0.008013"
Set-TextValue "E41" "4.44%"

Set-TextValue "E42" "6.81%"

Set-TextValue "D43" "0.007446"
Set-TextValue "E43" "1.07%"

Set-TextValue "D44" "0.008491"
Set-TextValue "E44" "11.60%"

Set-TextValue "D45" "0.3163"
Set-TextValue "E45" "4.68%"

Set-TextValue "D46" "0.00006635"
Set-TextValue "E46" "-0.93%"

Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.00%"

Set-TextValue "D48" "0.05491"
Set-TextValue "E48" "29.67%"

Set-TextValue "D49" "0.004005"
Set-TextValue "E49" "-4.66%"

Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "0.00%"

Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "0.00%"
